$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Names")

$ws.Range("A4").Value = "Hungary"
$ws.Range("B4").Value = "tu_cbfam_hu_HeadID"

$ws.Range("B2").Select()
